$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New headers for the straightness index columns
$ws.Range("N1").Value = "Straightness_Paired_Pre"
$ws.Range("O1").Value = "Straightness_Paired_Post"

# New data values (straightness index, pre/post) for rows 2-19
$straightness = @(
    @(0.0024821829999999998, 0.01205503),
    @(0.017634190000000001, 0.02527751),
    @(0.01138217, 0.0059585899999999997),
    @(0.015661370000000001, 0.0044013949999999998),
    @(0.018659889999999998, 0.016444529999999999),
    @(0.01051441, 0.0011427080000000001),
    @(0.0033264470000000002, 0.013293340000000001),
    @(0.01148972, 0.0079209039999999994),
    @(0.015130950000000001, 0.056988030000000002),
    @(0.0090667119999999993, 0.014699490000000001),
    @(0.0059884200000000004, 0.02223319),
    @(0.0084045969999999998, 0.0067090759999999996),
    @(0.0041845129999999999, 0.017935300000000001),
    @(0.0062429829999999997, 0.034996270000000003),
    @(0.004557635, 0.021763049999999999),
    @(0.032368510000000003, 0.056071160000000002),
    @(0.0058244580000000002, 0.01370646),
    @(0.016284449999999999, 0.010951610000000001)
)

$row = 2
foreach ($pair in $straightness) {
    $ws.Cells.Item($row, 14).Value = $pair[0]
    $ws.Cells.Item($row, 15).Value = $pair[1]
    $row++
}

# Update the selection to reflect where editing ended
$ws.Range("N22").Select()
